# Rename the inline logo pictures that live in the document's headers
# and footers:
#   - the Pearson Edexcel logo (footers)  : image1.png -> image2.png
#   - the BTEC logo               (headers): image2.jpg -> image1.jpg
#
# We walk every section's Headers/Footers collection (rather than
# assuming a fixed header1/footer1 layout) and use each picture's
# AlternativeText (the OOXML "descr" attribute, left untouched by this
# change) to decide which rename applies.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}
